# "Add files via upload" — turns the single-sheet "TEST CASE" workbook into a
# 3-sheet workbook: FULL CIRCUIT (renamed Sheet1), SERIAL TO PARALLEL, and
# DISPLAY BUFFER, each carrying its own truth-table data.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet, fix its selection -------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "FULL CIRCUIT"
[void]$ws1.Range("B2:G33").Select()

# --- Add "SERIAL TO PARALLEL" sheet ---------------------------------------
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "SERIAL TO PARALLEL"

# (COM ColumnWidth is in "characters"; the engine re-quantises to 1/6-character
# steps when serialising the stored OOXML `width`, so these inputs are chosen
# to land as close as possible to the original bestFit widths: 4, 6.109375,
# 9.44140625, 13.77734375, 8.)
$ws2.Columns.Item(2).ColumnWidth = 3.1666666666666665
$ws2.Columns.Item(3).ColumnWidth = 5.333333333333333
$ws2.Columns.Item(4).ColumnWidth = 8.666666666666666
$ws2.Columns.Item(5).ColumnWidth = 13
$ws2.Columns.Item(6).ColumnWidth = 7.166666666666667

$ws2.Range("B2").Value = "CLK"
$ws2.Range("C2").Value = "START"
$ws2.Range("D2").Value = "SERIAL_IN"
$ws2.Range("E2").Value = "PARALLEL_OUT"
$ws2.Range("F2").Value = "DISP_EN"

$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = "0xFCFC"
$ws2.Range("F3").Value = 1

$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = 0
$ws2.Range("D4").Value = 1
$ws2.Range("E4").Value = "0xFEFC"
$ws2.Range("F4").Value = 0

$ws2.Range("B5").Value = 0
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = "0xFEFC"
$ws2.Range("F5").Value = 0

$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = 0
$ws2.Range("D6").Value = 1
$ws2.Range("E6").Value = "0x1FC"
$ws2.Range("F6").Value = 0

$ws2.Range("B7").Value = 0
$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 1
$ws2.Range("E7").Value = "0x1FC"
$ws2.Range("F7").Value = 0

$ws2.Range("B8").Value = 1
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 0
$ws2.Range("E8").Value = "0x1FC"
$ws2.Range("F8").Value = 0

$ws2.Range("B9").Value = 0
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = "0x1FC"
$ws2.Range("F9").Value = 0

$ws2.Range("B10").Value = 1
$ws2.Range("C10").Value = 0
$ws2.Range("D10").Value = 1
$ws2.Range("E10").Value = "0xE0FC"
$ws2.Range("F10").Value = 0

$ws2.Range("B11").Value = 0
$ws2.Range("C11").Value = 0
$ws2.Range("D11").Value = 1
$ws2.Range("E11").Value = "0xE0FC"
$ws2.Range("F11").Value = 0

$ws2.Range("B12").Value = 1
$ws2.Range("C12").Value = 0
$ws2.Range("D12").Value = 0
$ws2.Range("E12").Value = "0x1FE"
$ws2.Range("F12").Value = 0

$ws2.Range("B13").Value = 0
$ws2.Range("C13").Value = 0
$ws2.Range("D13").Value = 0
$ws2.Range("E13").Value = "0x1FE"
$ws2.Range("F13").Value = 0

$ws2.Range("B14").Value = 1
$ws2.Range("C14").Value = 0
$ws2.Range("D14").Value = 0
$ws2.Range("E14").Value = "0xB601"
$ws2.Range("F14").Value = 0

$ws2.Range("B15").Value = 0
$ws2.Range("C15").Value = 0
$ws2.Range("D15").Value = 0
$ws2.Range("E15").Value = "0xB601"
$ws2.Range("F15").Value = 0

$ws2.Range("B16").Value = 1
$ws2.Range("C16").Value = 0
$ws2.Range("D16").Value = 1
$ws2.Range("E16").Value = "0xDA01"
$ws2.Range("F16").Value = 0

$ws2.Range("B17").Value = 0
$ws2.Range("C17").Value = 0
$ws2.Range("D17").Value = 1
$ws2.Range("E17").Value = "0xDA01"
$ws2.Range("F17").Value = 0

$ws2.Range("B18").Value = 1
$ws2.Range("C18").Value = 0
$ws2.Range("D18").Value = 0
$ws2.Range("E18").Value = "0xF6E0"
$ws2.Range("F18").Value = 1

[void]$ws2.Range("E20").Select()

# --- Add "DISPLAY BUFFER" sheet -------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "DISPLAY BUFFER"

$ws3.Range("B3").Value = "CLK"
$ws3.Range("C3").Value = "PARALLEL_IN"
$ws3.Range("D3").Value = "ENABLE"
$ws3.Range("E3").Value = "DISPLAY_OUT"

$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = "0xFCFC"
$ws3.Range("D4").Value = 1
$ws3.Range("E4").Value = "0xFCFC"

$ws3.Range("B5").Value = 0
$ws3.Range("C5").Value = "0xFCFC"
$ws3.Range("D5").Value = 1
$ws3.Range("E5").Value = "0xFCFC"

$ws3.Range("B6").Value = 1
$ws3.Range("C6").Value = "0xE0FC"
$ws3.Range("D6").Value = 1
$ws3.Range("E6").Value = "0xE0FC"

$ws3.Range("B7").Value = 0
$ws3.Range("C7").Value = "0xE0FC"
$ws3.Range("D7").Value = 1
$ws3.Range("E7").Value = "0xE0FC"

$ws3.Range("B8").Value = 1
$ws3.Range("C8").Value = "0xFEFC"
$ws3.Range("D8").Value = 1
$ws3.Range("E8").Value = "0xFEFC"

$ws3.Range("B9").Value = 0
$ws3.Range("C9").Value = "0xFEFC"
$ws3.Range("D9").Value = 1
$ws3.Range("E9").Value = "0xFEFC"

$ws3.Range("B10").Value = 1
$ws3.Range("C10").Value = "0xB660"
$ws3.Range("D10").Value = 1
$ws3.Range("E10").Value = "0xB660"

$ws3.Range("B11").Value = 0
$ws3.Range("C11").Value = "0xB660"
$ws3.Range("D11").Value = 1
$ws3.Range("E11").Value = "0xB660"

[void]$ws3.Range("B3:E11").Select()

# Leave the workbook on the original sheet/selection, matching the saved file.
[void]$ws1.Activate()
